$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "29.173.21"
$ws.Range("E2").Value = "  -2.13%  "

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.899.47"
$ws.Range("E3").Value = "  -2.72%  "

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.001"
$ws.Range("E4").Value = "  +0.08%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "331.92"
$ws.Range("E5").Value = "  -3.01%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.000"
$ws.Range("E6").Value = "  +0.07%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4606"
$ws.Range("E7").Value = "  -3.64%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.4136"
$ws.Range("E8").Value = "  -0.43%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "47.77"
$ws.Range("E9").Value = "  -1.53%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.08017"
$ws.Range("E10").Value = "  -3.05%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "1.012"
$ws.Range("E11").Value = "  -2.86%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "22.28"
$ws.Range("E12").Value = "  -1.86%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.894.69"
$ws.Range("E13").Value = "  -2.55%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.946"
$ws.Range("E14").Value = "  -4.26%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.122"
$ws.Range("E15").Value = "  -4.23%  "

# Row 16
$ws.Range("B16").Value = "Litecoin"
$ws.Range("C16").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "89.26"
$ws.Range("E16").Value = "  -3.47%  "

# Row 17
$ws.Range("B17").Value = "BinanceUSD"
$ws.Range("C17").Value = "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "1.000"
$ws.Range("E17").Value = "  -0.02%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.00001031"
$ws.Range("E18").Value = "  -3.19%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06577"
$ws.Range("E19").Value = "  -1.47%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "17.66"
$ws.Range("E20").Value = "  -2.58%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "1.002"
$ws.Range("E21").Value = "  +0.17%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "29.129.13"
$ws.Range("E22").Value = "  -2.13%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.502"
$ws.Range("E23").Value = "  -1.88%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "11.45"
$ws.Range("E24").Value = "  +1.47%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.198"
$ws.Range("E25").Value = "  -3.63%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.130.89"
$ws.Range("E26").Value = "  -1.80%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "156.37"
$ws.Range("E27").Value = "  -2.87%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "19.73"
$ws.Range("E28").Value = "  -2.39%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.126"
$ws.Range("E29").Value = "  -2.99%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "5.654"
$ws.Range("E30").Value = "  -0.91%  "

# Row 31
$ws.Range("E31").Value = "  -4.32%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.044"
$ws.Range("E32").Value = "  +2.41%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.09387"
$ws.Range("E33").Value = "  -2.56%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.418"
$ws.Range("E34").Value = "  -4.38%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "3.520"
$ws.Range("E35").Value = "  -4.41%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "5.358"
$ws.Range("E36").Value = "  -3.06%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.06084"
$ws.Range("E37").Value = "  -3.64%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.02240"
$ws.Range("E38").Value = "  -4.10%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "8.448"
$ws.Range("E39").Value = "  -1.16%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.178"
$ws.Range("E40").Value = "  -1.94%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.5843"
$ws.Range("E41").Value = "  -4.53%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.001"
$ws.Range("E42").Value = "  +0.21%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.1827"
$ws.Range("E43").Value = "  -3.84%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "10.16"
$ws.Range("E44").Value = "  -5.38%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.264"
$ws.Range("E45").Value = "  -2.10%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.07527"
$ws.Range("E46").Value = "  +1.82%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.310"
$ws.Range("E47").Value = "  -3.89%  "

# Row 48
$ws.Range("B48").Value = "Decentraland"
$ws.Range("C48").Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.5522"
$ws.Range("E48").Value = "  -3.53%  "

# Row 49
$ws.Range("B49").Value = "EnergySwap"
$ws.Range("C49").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "12.01"
$ws.Range("E49").Value = "  -4.78%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.935"
$ws.Range("E50").Value = "  -3.30%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "112.00"
$ws.Range("E51").Value = "  -2.17%  "
